# "link to survey deleted" — update the NO2 indicator worksheet inputs and
# refresh the sheet view (zoom + selection) to match the saved state, and
# turn on iterative calculation as the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three changed input cells -------------------------------
# D1: Cjm[NO2] measured value
$ws.Range("D1").Value = 5.13

# D14: Fmeteo constant (can be changed on Tygron)
$ws.Range("D14").Value = 1.8

# D18: Road length (m)
$ws.Range("D18").Value = 8003

# All dependent formulas (D3, D8, D9, D10, D11, D13, D15, D16, D17, D19)
# recalculate automatically from these inputs.

# --- Enable iterative calculation (calcPr iterate="1") ------------------
$excel.Iteration = $true

# --- Sheet view: zoom + final selected cell ------------------------------
$excel.ActiveWindow.Zoom = 145
$ws.Range("D19").Select() | Out-Null
